# Update the simulated power-flow voltage magnitude (vm_pu) results table
# on Sheet1 with the re-run values ("anadidas dos lineas dobles").
# Columns: A=index, B..H and L..O hold per-bus vm_pu values for rows 2-25.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 0.9773520444207801
$ws.Cells.Item(2, 4).Value = 0.9543152958162814
$ws.Cells.Item(2, 5).Value = 0.9582015157281115
$ws.Cells.Item(2, 7).Value = 0.9584562192782538
$ws.Cells.Item(2, 8).Value = 0.9551809252010338
$ws.Cells.Item(2, 12).Value = 0.9695248790727012
$ws.Cells.Item(2, 13).Value = 0.9338291813329535
$ws.Cells.Item(2, 14).Value = 0.9467395913474621
$ws.Cells.Item(2, 15).Value = 0.9504588672154154

$ws.Cells.Item(3, 3).Value = 0.9787449255614321
$ws.Cells.Item(3, 4).Value = 0.9566603145303347
$ws.Cells.Item(3, 5).Value = 0.9605777778798741
$ws.Cells.Item(3, 6).Value = 0.9999999999999997
$ws.Cells.Item(3, 7).Value = 0.9605458901460937
$ws.Cells.Item(3, 8).Value = 0.9575280710077952
$ws.Cells.Item(3, 12).Value = 0.9710155883320741
$ws.Cells.Item(3, 13).Value = 0.9363209389686052
$ws.Cells.Item(3, 14).Value = 0.9496437474603603
$ws.Cells.Item(3, 15).Value = 0.9526568489051052

$ws.Cells.Item(4, 3).Value = 0.9796948347367467
$ws.Cells.Item(4, 4).Value = 0.9583672781438106
$ws.Cells.Item(4, 5).Value = 0.9623419063884882
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.9621002288783299
$ws.Cells.Item(4, 8).Value = 0.9592365829542686
$ws.Cells.Item(4, 12).Value = 0.9720397813956198
$ws.Cells.Item(4, 13).Value = 0.9381380836259479
$ws.Cells.Item(4, 14).Value = 0.9518046843571263
$ws.Cells.Item(4, 15).Value = 0.9542938149014041

$ws.Cells.Item(5, 3).Value = 0.9801408591636842
$ws.Cells.Item(5, 4).Value = 0.9592141891442756
$ws.Cells.Item(5, 5).Value = 0.9632307717889166
$ws.Cells.Item(5, 6).Value = 0.9999999999999999
$ws.Cells.Item(5, 7).Value = 0.9628844653949439
$ws.Cells.Item(5, 8).Value = 0.9600842621610606
$ws.Cells.Item(5, 12).Value = 0.9725238731915415
$ws.Cells.Item(5, 13).Value = 0.9390410107474567
$ws.Cells.Item(5, 14).Value = 0.9528954492518774
$ws.Cells.Item(5, 15).Value = 0.9551205384500632

$ws.Cells.Item(6, 3).Value = 0.9803228908569025
$ws.Cells.Item(6, 4).Value = 0.9595709796282883
$ws.Cells.Item(6, 5).Value = 0.9636083939858976
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.9632178731693578
$ws.Cells.Item(6, 8).Value = 0.9604413762784975
$ws.Cells.Item(6, 12).Value = 0.9727222230810503
$ws.Cells.Item(6, 13).Value = 0.9394217159708234
$ws.Cells.Item(6, 14).Value = 0.953359304188821
$ws.Cells.Item(6, 15).Value = 0.9554721913041527

$ws.Cells.Item(7, 3).Value = 0.9803228908569025
$ws.Cells.Item(7, 4).Value = 0.9595709796282883
$ws.Cells.Item(7, 5).Value = 0.9636083939858976
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.9632178731693578
$ws.Cells.Item(7, 8).Value = 0.9604413762784975
$ws.Cells.Item(7, 12).Value = 0.9727222230810503
$ws.Cells.Item(7, 13).Value = 0.9394217159708234
$ws.Cells.Item(7, 14).Value = 0.953359304188821
$ws.Cells.Item(7, 15).Value = 0.9554721913041527

$ws.Cells.Item(8, 3).Value = 0.9802785576336249
$ws.Cells.Item(8, 4).Value = 0.9594834102319267
$ws.Cells.Item(8, 5).Value = 0.9635155267855273
$ws.Cells.Item(8, 6).Value = 0.9999999999999999
$ws.Cells.Item(8, 7).Value = 0.9631358659848587
$ws.Cells.Item(8, 8).Value = 0.9603537274506913
$ws.Cells.Item(8, 12).Value = 0.972673868324171
$ws.Cells.Item(8, 13).Value = 0.939328258457981
$ws.Cells.Item(8, 14).Value = 0.953245203260354
$ws.Cells.Item(8, 15).Value = 0.9553856858552044

$ws.Cells.Item(9, 3).Value = 0.9795875025155213
$ws.Cells.Item(9, 4).Value = 0.9581684574427475
$ws.Cells.Item(9, 5).Value = 0.9621346472824567
$ws.Cells.Item(9, 7).Value = 0.9619174739392945
$ws.Cells.Item(9, 8).Value = 0.9590375819092017
$ws.Cells.Item(9, 12).Value = 0.9719236381435616
$ws.Cells.Item(9, 13).Value = 0.937926253259975
$ws.Cells.Item(9, 14).Value = 0.9515505507268847
$ws.Cells.Item(9, 15).Value = 0.9541012404150635

$ws.Cells.Item(10, 3).Value = 0.9774625065754827
$ws.Cells.Item(10, 4).Value = 0.9544961363254286
$ws.Cells.Item(10, 5).Value = 0.9583831256041925
$ws.Cells.Item(10, 6).Value = 0.9999999999999999
$ws.Cells.Item(10, 7).Value = 0.9586157760056971
$ws.Cells.Item(10, 8).Value = 0.9553619297449185
$ws.Cells.Item(10, 12).Value = 0.9696427387633458
$ws.Cells.Item(10, 13).Value = 0.9340211795739423
$ws.Cells.Item(10, 14).Value = 0.9469613200609243
$ws.Cells.Item(10, 15).Value = 0.950626595704821

$ws.Cells.Item(11, 3).Value = 0.9748689292287538
$ws.Cells.Item(11, 4).Value = 0.9504114862214637
$ws.Cells.Item(11, 5).Value = 0.9543341647077224
$ws.Cells.Item(11, 6).Value = 0.9999999999999999
$ws.Cells.Item(11, 7).Value = 0.9550638063553704
$ws.Cells.Item(11, 8).Value = 0.9512735745833356
$ws.Cells.Item(11, 12).Value = 0.9668868274139121
$ws.Cells.Item(11, 13).Value = 0.9296895161211729
$ws.Cells.Item(11, 14).Value = 0.9420249381170074
$ws.Cells.Item(11, 15).Value = 0.9468959243232355

$ws.Cells.Item(12, 3).Value = 0.9734392118652518
$ws.Cells.Item(12, 4).Value = 0.94826691944806
$ws.Cells.Item(12, 5).Value = 0.9522449671963078
$ws.Cells.Item(12, 6).Value = 0.9999999999999999
$ws.Cells.Item(12, 7).Value = 0.9532350681961522
$ws.Cells.Item(12, 8).Value = 0.9491270625408735
$ws.Cells.Item(12, 12).Value = 0.9653751663401665
$ws.Cells.Item(12, 13).Value = 0.9274186317852324
$ws.Cells.Item(12, 14).Value = 0.939482498783516
$ws.Cells.Item(12, 15).Value = 0.944977434138223

$ws.Cells.Item(13, 3).Value = 0.9722409428425127
$ws.Cells.Item(13, 4).Value = 0.946510755333466
$ws.Cells.Item(13, 5).Value = 0.9505489586392033
$ws.Cells.Item(13, 7).Value = 0.9517523369277658
$ws.Cells.Item(13, 8).Value = 0.9473693054650555
$ws.Cells.Item(13, 12).Value = 0.9641111258541468
$ws.Cells.Item(13, 13).Value = 0.9255603389508904
$ws.Cells.Item(13, 14).Value = 0.9374203080226988
$ws.Cells.Item(13, 15).Value = 0.9434228634192503

$ws.Cells.Item(14, 3).Value = 0.9720230241673952
$ws.Cells.Item(14, 4).Value = 0.9461948207369479
$ws.Cells.Item(14, 5).Value = 0.9502451127247138
$ws.Cells.Item(14, 7).Value = 0.951486867414099
$ws.Cells.Item(14, 8).Value = 0.9470530842942034
$ws.Cells.Item(14, 12).Value = 0.9638814887877987
$ws.Cells.Item(14, 13).Value = 0.925226140979285
$ws.Cells.Item(14, 14).Value = 0.9370510037041635
$ws.Cells.Item(14, 15).Value = 0.9431446122253315

$ws.Cells.Item(15, 3).Value = 0.9723363685275841
$ws.Cells.Item(15, 4).Value = 0.946649416628798
$ws.Cells.Item(15, 5).Value = 0.9506824313911408
$ws.Cells.Item(15, 6).Value = 0.9999999999999999
$ws.Cells.Item(15, 7).Value = 0.9518689675980906
$ws.Cells.Item(15, 8).Value = 0.9475080925356868
$ws.Cells.Item(15, 12).Value = 0.964211705185073
$ws.Cells.Item(15, 13).Value = 0.9257070259434191
$ws.Cells.Item(15, 14).Value = 0.9375825483364657
$ws.Cells.Item(15, 15).Value = 0.9435451170423201

$ws.Cells.Item(16, 3).Value = 0.9723333995106044
$ws.Cells.Item(16, 4).Value = 0.9466450994696032
$ws.Cells.Item(16, 5).Value = 0.9506782746857758
$ws.Cells.Item(16, 6).Value = 0.9999999999999999
$ws.Cells.Item(16, 7).Value = 0.9518653352548426
$ws.Cells.Item(16, 8).Value = 0.9475037714605324
$ws.Cells.Item(16, 12).Value = 0.9642085756135909
$ws.Cells.Item(16, 13).Value = 0.9257024588139902
$ws.Cells.Item(16, 14).Value = 0.9375774956079905
$ws.Cells.Item(16, 15).Value = 0.9435413095084139

$ws.Cells.Item(17, 3).Value = 0.972314438727412
$ws.Cells.Item(17, 4).Value = 0.9466175336224966
$ws.Cells.Item(17, 5).Value = 0.950651735020445
$ws.Cells.Item(17, 7).Value = 0.9518421437501032
$ws.Cells.Item(17, 8).Value = 0.9474761806093129
$ws.Cells.Item(17, 12).Value = 0.9641885898095401
$ws.Cells.Item(17, 13).Value = 0.9256732970038003
$ws.Cells.Item(17, 14).Value = 0.9375452352155481
$ws.Cells.Item(17, 15).Value = 0.943516999565488

$ws.Cells.Item(18, 3).Value = 0.9732968422111016
$ws.Cells.Item(18, 4).Value = 0.948056471168108
$ws.Cells.Item(18, 5).Value = 0.9520410679423061
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 0.9530567247482383
$ws.Cells.Item(18, 8).Value = 0.9489164233699072
$ws.Cells.Item(18, 12).Value = 0.96522485571338
$ws.Cells.Item(18, 13).Value = 0.9271958873455026
$ws.Cells.Item(18, 14).Value = 0.9392344991384203
$ws.Cells.Item(18, 15).Value = 0.9447904079111434

$ws.Cells.Item(19, 3).Value = 0.9738791038624653
$ws.Cells.Item(19, 4).Value = 0.9489205038430967
$ws.Cells.Item(19, 5).Value = 0.9528794316072505
$ws.Cells.Item(19, 7).Value = 0.9537901638074934
$ws.Cells.Item(19, 8).Value = 0.9497812397817553
$ws.Cells.Item(19, 12).Value = 0.9658398300098027
$ws.Cells.Item(19, 13).Value = 0.9281105119794807
$ws.Cells.Item(19, 14).Value = 0.9402543339191191
$ws.Cells.Item(19, 15).Value = 0.9455596320780167

$ws.Cells.Item(20, 3).Value = 0.9738748813254874
$ws.Cells.Item(20, 4).Value = 0.9489142052336506
$ws.Cells.Item(20, 5).Value = 0.9528733082505921
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 0.9537848053471792
$ws.Cells.Item(20, 8).Value = 0.9497749354590389
$ws.Cells.Item(20, 12).Value = 0.9658353679381856
$ws.Cells.Item(20, 13).Value = 0.9281038435118653
$ws.Cells.Item(20, 14).Value = 0.9402468836844801
$ws.Cells.Item(20, 15).Value = 0.9455540114247609

$ws.Cells.Item(21, 3).Value = 0.9723126027603244
$ws.Cells.Item(21, 4).Value = 0.9466148648389406
$ws.Cells.Item(21, 5).Value = 0.9506491657383351
$ws.Cells.Item(21, 6).Value = 0.9999999999999998
$ws.Cells.Item(21, 7).Value = 0.9518398986208889
$ws.Cells.Item(21, 8).Value = 0.9474735094049876
$ws.Cells.Item(21, 12).Value = 0.9641866546188596
$ws.Cells.Item(21, 13).Value = 0.9256704737203618
$ws.Cells.Item(21, 14).Value = 0.9375421121319083
$ws.Cells.Item(21, 15).Value = 0.9435146461719072

$ws.Cells.Item(22, 3).Value = 0.9706597409360523
$ws.Cells.Item(22, 4).Value = 0.9442393153608207
$ws.Cells.Item(22, 5).Value = 0.9483722267301928
$ws.Cells.Item(22, 7).Value = 0.9498516146099044
$ws.Cells.Item(22, 8).Value = 0.945095805140664
$ws.Cells.Item(22, 12).Value = 0.9624463749959533
$ws.Cells.Item(22, 13).Value = 0.9231582522576942
$ws.Cells.Item(22, 14).Value = 0.934775485679737
$ws.Cells.Item(22, 15).Value = 0.9414311248838643

$ws.Cells.Item(23, 3).Value = 0.9693014104230095
$ws.Cells.Item(23, 4).Value = 0.9423224331493186
$ws.Cells.Item(23, 5).Value = 0.9465481844322358
$ws.Cells.Item(23, 7).Value = 0.9482607458305827
$ws.Cells.Item(23, 8).Value = 0.9431771841856078
$ws.Cells.Item(23, 12).Value = 0.9610186996633683
$ws.Cells.Item(23, 13).Value = 0.9211321842097246
$ws.Cells.Item(23, 14).Value = 0.9325605384440879
$ws.Cells.Item(23, 15).Value = 0.9397649083974079

$ws.Cells.Item(24, 3).Value = 0.968775405443835
$ws.Cells.Item(24, 4).Value = 0.9415874928651377
$ws.Cells.Item(24, 5).Value = 0.9458516534909652
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 0.9476536932431231
$ws.Cells.Item(24, 8).Value = 0.9424415772602138
$ws.Cells.Item(24, 12).Value = 0.9604663610606421
$ws.Cells.Item(24, 13).Value = 0.9203556068545187
$ws.Cells.Item(24, 14).Value = 0.9317150098529556
$ws.Cells.Item(24, 15).Value = 0.939129287291356

$ws.Cells.Item(25, 3).Value = 0.971402512821942
$ws.Cells.Item(25, 4).Value = 0.9453004487010482
$ws.Cells.Item(25, 5).Value = 0.9493869089884963
$ws.Cells.Item(25, 7).Value = 0.9507373239545897
$ws.Cells.Item(25, 8).Value = 0.9461579010015644
$ws.Cells.Item(25, 12).Value = 0.9632279792700686
$ws.Cells.Item(25, 13).Value = 0.9242802332416815
$ws.Cells.Item(25, 14).Value = 0.9360081295287115
$ws.Cells.Item(25, 15).Value = 0.9423591046471748
